$d = $word.ActiveDocument

# Locate the anchor paragraph: "sudo update-alternatives --set php /usr/bin/php8.3"
# (the first/only occurrence, found right after the "Command Line:" block).
$rng = $d.Content
$found = $rng.Find.Execute("sudo update-alternatives --set php /usr/bin/php8.3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$lines = @(
    "",
    "/*",
    "//for enable php 7.4",
    "sudo update-alternatives --list php",
    "sudo a2dismod php8.3",
    "sudo a2enmod php7.4",
    "sudo update-alternatives --set php /usr/bin/php7.4",
    "systemctl restart apache2",
    "",
    "",
    "//for enable php 8.3",
    "sudo update-alternatives --list php",
    "sudo a2dismod php7.4",
    "sudo a2enmod php8.3",
    "sudo update-alternatives --set php /usr/bin/php8.3",
    "systemctl restart apache2"
)

foreach ($line in $lines) {
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null
    if ($line -ne "") {
        $rng.InsertAfter($line)
        $rng.Collapse(0)
    }
}

# Flip the Normal style's OverflowPunct (paragraph format) false -> true.
$style = $d.Styles("Normal")
$style.ParagraphFormat.HangingPunctuation = $true
